# Updated exporters to support context data
# Appends 4 new order lines (Natalie's juices) to the bottom of the
# Performance Food / Downtown order sheet, extending the used range
# from A1:E16 to A1:E20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append, in sheet (SKU, Name, Quantity, Cost Per, Total Cost) order.
$newRows = @(
    @("TN330", "Natalie's - Honey Tangerine", "1", "14.57", "14.57"),
    @("AH252", "Natalie's - Orange Juice",    "2", "24.50", "49.00"),
    @("TN454", "Natalie's - Orange Mango",    "1", "13.38", "13.38"),
    @("TN362", "Natalie's - Orange Pineapple","1", "13.38", "13.38")
)

$startRow = 17
$endRow = $startRow + $newRows.Length - 1

# The source data stores every column (including Quantity/Cost Per/Total
# Cost) as text, not numbers. Force the destination range to a text
# number format before writing so the quantity/price strings aren't
# auto-coerced into numeric cells, then strip the format back off so the
# cells end up with the workbook's default (unstyled) formatting, same
# as every other row already on the sheet.
$target = $ws.Range("A$startRow`:E$endRow")
$target.NumberFormat = "@"

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    $ws.Range("A$r").Value = $rowData[0]
    $ws.Range("B$r").Value = $rowData[1]
    $ws.Range("C$r").Value = $rowData[2]
    $ws.Range("D$r").Value = $rowData[3]
    $ws.Range("E$r").Value = $rowData[4]
}

$target.ClearFormats()
